$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 0.291
$ws.Cells.Item(4, 5).Value = 0.158
$ws.Cells.Item(4, 8).Value = 0.193
$ws.Cells.Item(4, 11).Value = 0.346
$ws.Cells.Item(4, 12).Value = 0.098
$ws.Cells.Item(4, 13).Value = 0.313
$ws.Cells.Item(4, 14).Value = 0.266
$ws.Cells.Item(4, 15).Value = 0.021
$ws.Cells.Item(4, 16).Value = 0.146
$ws.Cells.Item(4, 17).Value = 0.538
$ws.Cells.Item(4, 18).Value = 0.217
$ws.Cells.Item(4, 19).Value = 0.466
$ws.Cells.Item(4, 20).Value = 0.275
$ws.Cells.Item(4, 23).Value = 0.251
$ws.Cells.Item(4, 24).Value = 0.042
$ws.Cells.Item(4, 25).Value = 0.204
$ws.Cells.Item(4, 26).Value = 0.445
$ws.Cells.Item(4, 27).Value = 0.128
$ws.Cells.Item(4, 28).Value = 0.358
$ws.Cells.Item(4, 29).Value = 0.127
$ws.Cells.Item(4, 31).Value = 0.077
$ws.Cells.Item(4, 32).Value = 0.721
$ws.Cells.Item(4, 35).Value = 0.675
$ws.Cells.Item(4, 36).Value = 0.157
$ws.Cells.Item(4, 37).Value = 0.396
$ws.Cells.Item(4, 38).Value = 0.697
$ws.Cells.Item(4, 39).Value = 0.117
$ws.Cells.Item(4, 40).Value = 0.342
$ws.Cells.Item(4, 41).Value = 0.698
# Row 5
$ws.Cells.Item(5, 2).Value = 0.825
$ws.Cells.Item(5, 3).Value = 0.144
$ws.Cells.Item(5, 4).Value = 0.38
$ws.Cells.Item(5, 5).Value = 0.7
$ws.Cells.Item(5, 6).Value = 0.21
$ws.Cells.Item(5, 7).Value = 0.458
$ws.Cells.Item(5, 8).Value = 0.85
$ws.Cells.Item(5, 9).Value = 0.128
$ws.Cells.Item(5, 10).Value = 0.357
$ws.Cells.Item(5, 11).Value = 0.675
$ws.Cells.Item(5, 12).Value = 0.219
$ws.Cells.Item(5, 13).Value = 0.468
$ws.Cells.Item(5, 14).Value = 0.825
$ws.Cells.Item(5, 15).Value = 0.144
$ws.Cells.Item(5, 16).Value = 0.38
$ws.Cells.Item(5, 17).Value = 0.6
$ws.Cells.Item(5, 18).Value = 0.24
$ws.Cells.Item(5, 19).Value = 0.49
$ws.Cells.Item(5, 20).Value = 0.575
$ws.Cells.Item(5, 21).Value = 0.244
$ws.Cells.Item(5, 22).Value = 0.494
$ws.Cells.Item(5, 23).Value = 0.775
$ws.Cells.Item(5, 24).Value = 0.174
$ws.Cells.Item(5, 25).Value = 0.418
$ws.Cells.Item(5, 26).Value = 0.825
$ws.Cells.Item(5, 27).Value = 0.144
$ws.Cells.Item(5, 28).Value = 0.38
$ws.Cells.Item(5, 29).Value = 0.775
$ws.Cells.Item(5, 30).Value = 0.174
$ws.Cells.Item(5, 31).Value = 0.418
$ws.Cells.Item(5, 32).Value = 0.975
$ws.Cells.Item(5, 33).Value = 0.024
$ws.Cells.Item(5, 34).Value = 0.156
$ws.Cells.Item(5, 35).Value = 0.8
$ws.Cells.Item(5, 36).Value = 0.16
$ws.Cells.Item(5, 37).Value = 0.4
$ws.Cells.Item(5, 38).Value = 0.925
$ws.Cells.Item(5, 39).Value = 0.06900000000000001
$ws.Cells.Item(5, 40).Value = 0.263
$ws.Cells.Item(5, 41).Value = 0.9
# Row 6
$ws.Cells.Item(6, 2).Value = 0.43
$ws.Cells.Item(6, 5).Value = 0.258
$ws.Cells.Item(6, 8).Value = 0.315
$ws.Cells.Item(6, 11).Value = 0.457
$ws.Cells.Item(6, 14).Value = 0.402
$ws.Cells.Item(6, 17).Value = 0.5669999999999999
$ws.Cells.Item(6, 20).Value = 0.372
$ws.Cells.Item(6, 23).Value = 0.379
$ws.Cells.Item(6, 26).Value = 0.578
$ws.Cells.Item(6, 29).Value = 0.218
$ws.Cells.Item(6, 32).Value = 0.829
$ws.Cells.Item(6, 35).Value = 0.732
$ws.Cells.Item(6, 38).Value = 0.795
$ws.Cells.Item(6, 41).Value = 0.785
# Row 7
$ws.Cells.Item(7, 2).Value = 0.604
$ws.Cells.Item(7, 5).Value = 0.415
$ws.Cells.Item(7, 8).Value = 0.506
$ws.Cells.Item(7, 11).Value = 0.5669999999999999
$ws.Cells.Item(7, 14).Value = 0.581
$ws.Cells.Item(7, 17).Value = 0.586
$ws.Cells.Item(7, 20).Value = 0.472
$ws.Cells.Item(7, 23).Value = 0.547
$ws.Cells.Item(7, 26).Value = 0.705
$ws.Cells.Item(7, 29).Value = 0.384
$ws.Cells.Item(7, 32).Value = 0.911
$ws.Cells.Item(7, 35).Value = 0.771
$ws.Cells.Item(7, 38).Value = 0.868
$ws.Cells.Item(7, 41).Value = 0.85
# Row 8
$ws.Cells.Item(8, 2).Value = 0.744
$ws.Cells.Item(8, 3).Value = 0.148
$ws.Cells.Item(8, 4).Value = 0.385
$ws.Cells.Item(8, 5).Value = 0.592
$ws.Cells.Item(8, 7).Value = 0.432
$ws.Cells.Item(8, 8).Value = 0.743
$ws.Cells.Item(8, 9).Value = 0.137
$ws.Cells.Item(8, 10).Value = 0.37
$ws.Cells.Item(8, 11).Value = 0.597
$ws.Cells.Item(8, 12).Value = 0.198
$ws.Cells.Item(8, 13).Value = 0.445
$ws.Cells.Item(8, 14).Value = 0.749
$ws.Cells.Item(8, 15).Value = 0.146
$ws.Cells.Item(8, 16).Value = 0.382
$ws.Cells.Item(8, 17).Value = 0.572
$ws.Cells.Item(8, 18).Value = 0.227
$ws.Cells.Item(8, 19).Value = 0.477
$ws.Cells.Item(8, 20).Value = 0.501
$ws.Cells.Item(8, 21).Value = 0.208
$ws.Cells.Item(8, 22).Value = 0.456
$ws.Cells.Item(8, 23).Value = 0.695
$ws.Cells.Item(8, 24).Value = 0.165
$ws.Cells.Item(8, 25).Value = 0.406
$ws.Cells.Item(8, 26).Value = 0.763
$ws.Cells.Item(8, 27).Value = 0.145
$ws.Cells.Item(8, 28).Value = 0.381
$ws.Cells.Item(8, 29).Value = 0.645
$ws.Cells.Item(8, 30).Value = 0.171
$ws.Cells.Item(8, 31).Value = 0.413
$ws.Cells.Item(8, 32).Value = 0.898
$ws.Cells.Item(8, 33).Value = 0.045
$ws.Cells.Item(8, 34).Value = 0.211
$ws.Cells.Item(8, 35).Value = 0.791
$ws.Cells.Item(8, 36).Value = 0.16
$ws.Cells.Item(8, 37).Value = 0.4
$ws.Cells.Item(8, 38).Value = 0.888
$ws.Cells.Item(8, 39).Value = 0.076
$ws.Cells.Item(8, 40).Value = 0.276
$ws.Cells.Item(8, 41).Value = 0.859
# Row 9
$ws.Cells.Item(9, 2).Value = 0.65
$ws.Cells.Item(9, 3).Value = 0.227
$ws.Cells.Item(9, 4).Value = 0.477
$ws.Cells.Item(9, 5).Value = 0.475
$ws.Cells.Item(9, 8).Value = 0.625
$ws.Cells.Item(9, 9).Value = 0.234
$ws.Cells.Item(9, 10).Value = 0.484
$ws.Cells.Item(9, 11).Value = 0.5
$ws.Cells.Item(9, 14).Value = 0.65
$ws.Cells.Item(9, 15).Value = 0.227
$ws.Cells.Item(9, 16).Value = 0.477
$ws.Cells.Item(9, 17).Value = 0.525
$ws.Cells.Item(9, 20).Value = 0.4
$ws.Cells.Item(9, 21).Value = 0.24
$ws.Cells.Item(9, 22).Value = 0.49
$ws.Cells.Item(9, 23).Value = 0.575
$ws.Cells.Item(9, 24).Value = 0.244
$ws.Cells.Item(9, 25).Value = 0.494
$ws.Cells.Item(9, 26).Value = 0.675
$ws.Cells.Item(9, 27).Value = 0.219
$ws.Cells.Item(9, 28).Value = 0.468
$ws.Cells.Item(9, 29).Value = 0.525
$ws.Cells.Item(9, 30).Value = 0.249
$ws.Cells.Item(9, 31).Value = 0.499
$ws.Cells.Item(9, 32).Value = 0.775
$ws.Cells.Item(9, 33).Value = 0.174
$ws.Cells.Item(9, 34).Value = 0.418
$ws.Cells.Item(9, 35).Value = 0.775
$ws.Cells.Item(9, 36).Value = 0.174
$ws.Cells.Item(9, 37).Value = 0.418
$ws.Cells.Item(9, 38).Value = 0.825
$ws.Cells.Item(9, 39).Value = 0.144
$ws.Cells.Item(9, 40).Value = 0.38
$ws.Cells.Item(9, 41).Value = 0.792
# Row 10
$ws.Cells.Item(10, 2).Value = 0.775
$ws.Cells.Item(10, 3).Value = 0.174
$ws.Cells.Item(10, 4).Value = 0.418
$ws.Cells.Item(10, 5).Value = 0.625
$ws.Cells.Item(10, 6).Value = 0.234
$ws.Cells.Item(10, 7).Value = 0.484
$ws.Cells.Item(10, 8).Value = 0.775
$ws.Cells.Item(10, 9).Value = 0.174
$ws.Cells.Item(10, 10).Value = 0.418
$ws.Cells.Item(10, 11).Value = 0.675
$ws.Cells.Item(10, 12).Value = 0.219
$ws.Cells.Item(10, 13).Value = 0.468
$ws.Cells.Item(10, 14).Value = 0.8
$ws.Cells.Item(10, 15).Value = 0.16
$ws.Cells.Item(10, 16).Value = 0.4
$ws.Cells.Item(10, 17).Value = 0.6
$ws.Cells.Item(10, 18).Value = 0.24
$ws.Cells.Item(10, 19).Value = 0.49
$ws.Cells.Item(10, 20).Value = 0.575
$ws.Cells.Item(10, 21).Value = 0.244
$ws.Cells.Item(10, 22).Value = 0.494
$ws.Cells.Item(10, 23).Value = 0.775
$ws.Cells.Item(10, 24).Value = 0.174
$ws.Cells.Item(10, 25).Value = 0.418
$ws.Cells.Item(10, 26).Value = 0.825
$ws.Cells.Item(10, 27).Value = 0.144
$ws.Cells.Item(10, 28).Value = 0.38
$ws.Cells.Item(10, 29).Value = 0.65
$ws.Cells.Item(10, 30).Value = 0.227
$ws.Cells.Item(10, 31).Value = 0.477
$ws.Cells.Item(10, 32).Value = 0.975
$ws.Cells.Item(10, 33).Value = 0.024
$ws.Cells.Item(10, 34).Value = 0.156
$ws.Cells.Item(10, 35).Value = 0.8
$ws.Cells.Item(10, 36).Value = 0.16
$ws.Cells.Item(10, 37).Value = 0.4
$ws.Cells.Item(10, 38).Value = 0.925
$ws.Cells.Item(10, 39).Value = 0.06900000000000001
$ws.Cells.Item(10, 40).Value = 0.263
$ws.Cells.Item(10, 41).Value = 0.9
# Row 11
$ws.Cells.Item(11, 2).Value = 0.825
$ws.Cells.Item(11, 3).Value = 0.144
$ws.Cells.Item(11, 4).Value = 0.38
$ws.Cells.Item(11, 5).Value = 0.7
$ws.Cells.Item(11, 6).Value = 0.21
$ws.Cells.Item(11, 7).Value = 0.458
$ws.Cells.Item(11, 8).Value = 0.85
$ws.Cells.Item(11, 9).Value = 0.128
$ws.Cells.Item(11, 10).Value = 0.357
$ws.Cells.Item(11, 11).Value = 0.675
$ws.Cells.Item(11, 12).Value = 0.219
$ws.Cells.Item(11, 13).Value = 0.468
$ws.Cells.Item(11, 14).Value = 0.825
$ws.Cells.Item(11, 15).Value = 0.144
$ws.Cells.Item(11, 16).Value = 0.38
$ws.Cells.Item(11, 17).Value = 0.6
$ws.Cells.Item(11, 18).Value = 0.24
$ws.Cells.Item(11, 19).Value = 0.49
$ws.Cells.Item(11, 20).Value = 0.575
$ws.Cells.Item(11, 21).Value = 0.244
$ws.Cells.Item(11, 22).Value = 0.494
$ws.Cells.Item(11, 23).Value = 0.775
$ws.Cells.Item(11, 24).Value = 0.174
$ws.Cells.Item(11, 25).Value = 0.418
$ws.Cells.Item(11, 26).Value = 0.825
$ws.Cells.Item(11, 27).Value = 0.144
$ws.Cells.Item(11, 28).Value = 0.38
$ws.Cells.Item(11, 29).Value = 0.7
$ws.Cells.Item(11, 30).Value = 0.21
$ws.Cells.Item(11, 31).Value = 0.458
$ws.Cells.Item(11, 32).Value = 0.975
$ws.Cells.Item(11, 33).Value = 0.024
$ws.Cells.Item(11, 34).Value = 0.156
$ws.Cells.Item(11, 35).Value = 0.8
$ws.Cells.Item(11, 36).Value = 0.16
$ws.Cells.Item(11, 37).Value = 0.4
$ws.Cells.Item(11, 38).Value = 0.925
$ws.Cells.Item(11, 39).Value = 0.06900000000000001
$ws.Cells.Item(11, 40).Value = 0.263
$ws.Cells.Item(11, 41).Value = 0.9
# Row 12
$ws.Cells.Item(12, 2).Value = 1.394
$ws.Cells.Item(12, 3).Value = 0.724
$ws.Cells.Item(12, 4).Value = 0.851
$ws.Cells.Item(12, 5).Value = 1.643
$ws.Cells.Item(12, 6).Value = 1.087
$ws.Cells.Item(12, 7).Value = 1.042
$ws.Cells.Item(12, 8).Value = 1.559
$ws.Cells.Item(12, 9).Value = 1.247
$ws.Cells.Item(12, 10).Value = 1.116
$ws.Cells.Item(12, 11).Value = 1.407
$ws.Cells.Item(12, 12).Value = 0.538
$ws.Cells.Item(12, 13).Value = 0.733
$ws.Cells.Item(12, 14).Value = 1.333
$ws.Cells.Item(12, 15).Value = 0.525
$ws.Cells.Item(12, 16).Value = 0.725
$ws.Cells.Item(12, 26).Value = 1.242
$ws.Cells.Item(12, 27).Value = 0.305
$ws.Cells.Item(12, 28).Value = 0.552
$ws.Cells.Item(12, 29).Value = 2.032
$ws.Cells.Item(12, 30).Value = 3.902
$ws.Cells.Item(12, 31).Value = 1.975
$ws.Cells.Item(12, 32).Value = 1.231
$ws.Cells.Item(12, 33).Value = 0.229
$ws.Cells.Item(12, 34).Value = 0.478
$ws.Cells.Item(12, 35).Value = 1.031
$ws.Cells.Item(12, 36).Value = 0.03
$ws.Cells.Item(12, 37).Value = 0.174
$ws.Cells.Item(12, 38).Value = 1.108
$ws.Cells.Item(12, 39).Value = 0.096
$ws.Cells.Item(12, 40).Value = 0.311
$ws.Cells.Item(12, 41).Value = 1.123
# Row 13
$ws.Cells.Item(13, 2).Value = 3.525
$ws.Cells.Item(13, 3).Value = 1.399
$ws.Cells.Item(13, 4).Value = 1.183
$ws.Cells.Item(13, 5).Value = 4.647
$ws.Cells.Item(13, 6).Value = 0.405
$ws.Cells.Item(13, 8).Value = 4.595
$ws.Cells.Item(13, 9).Value = 0.673
$ws.Cells.Item(13, 10).Value = 0.821
$ws.Cells.Item(13, 11).Value = 2.3
$ws.Cells.Item(13, 12).Value = 0.61
$ws.Cells.Item(13, 13).Value = 0.781
$ws.Cells.Item(13, 14).Value = 3.275
$ws.Cells.Item(13, 15).Value = 0.749
$ws.Cells.Item(13, 16).Value = 0.866
$ws.Cells.Item(13, 26).Value = 2.846
$ws.Cells.Item(13, 27).Value = 4.079
$ws.Cells.Item(13, 28).Value = 2.02
$ws.Cells.Item(13, 29).Value = 6.375
$ws.Cells.Item(13, 30).Value = 2.334
$ws.Cells.Item(13, 31).Value = 1.528
$ws.Cells.Item(13, 32).Value = 1.65
$ws.Cells.Item(13, 33).Value = 0.727
$ws.Cells.Item(13, 34).Value = 0.853
$ws.Cells.Item(13, 35).Value = 1.25
$ws.Cells.Item(13, 36).Value = 0.188
$ws.Cells.Item(13, 37).Value = 0.433
$ws.Cells.Item(13, 38).Value = 1.65
$ws.Cells.Item(13, 39).Value = 0.827
$ws.Cells.Item(13, 40).Value = 0.91
$ws.Cells.Item(13, 41).Value = 1.517
